$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "heuristic3"

$values = @(2,3,3,4,6,9,8,9,8,9,11,11,11,11,11,15,15,13,14,14)
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(2, $col).Value = $values[$i]
}
